$d = $word.ActiveDocument

# 1. Merge the three runs of the title paragraph ("Activity: Problem Sol" /
#    bookmark _GoBack / "ving " / "10/1/13") into a single run reading
#    "Activity: Problem Solving 10/1/13". Because the existing _GoBack
#    bookmark sits inside the text being replaced, Find/Replace removes it
#    for us as part of collapsing the runs.
$d.Content.Find.Execute("Problem Solving 10/1/13", $true, $false, $false, $false, $false, $true, 1, $false, "Problem Solving 10/1/13", 2) | Out-Null

# 2. Insert a new paragraph right after the title paragraph containing the
#    activity's name, with no special paragraph style/alignment.
$titlePara = $d.Paragraphs(2)
$titlePara.Range.InsertParagraphAfter()
$subtitlePara = $d.Paragraphs(3)
$subtitlePara.Range.Style = "Normal"
$subtitlePara.Range.ParagraphFormat.Alignment = 0
$subtitlePara.Range.Text = "A Cat, a Parrot, and a Bag of Seed"

# 3. Paragraph 4 is the pre-existing blank paragraph; leave it untouched.

# 4. Insert the explanatory paragraph right after the blank paragraph.
$blankPara = $d.Paragraphs(4)
$blankPara.Range.InsertParagraphAfter()
$bodyPara = $d.Paragraphs(5)
$bodyPara.Range.Style = "Normal"
$bodyPara.Range.ParagraphFormat.Alignment = 0
$bodyPara.Range.Text = "A man must transport items A, B, and C individually where neither A and B, nor B and C may left alone with each other. Items A, B, and C are a cat, a parrot, and a bag of seed, respectively. It is not stated in the problem but I can infer that the cat and bag of seeds pose no threat to each other because cats don’t eat seeds and seeds don’t eat cats. So, the man should not be constrained by leaving items A and C alone with each other."

# 5. Append a final, otherwise-empty paragraph that only carries the
#    _GoBack bookmark (mirrors where Word leaves it after the last edit).
#    Placing a zero-length bookmark at the very end of the document tends
#    to land the bookmark in the wrong spot, so first insert a couple of
#    placeholder characters, plant the bookmark at the start of them, then
#    delete the placeholder text again.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$bookmarkParaIndex = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs($bookmarkParaIndex)
$bookmarkPara.Range.InsertBefore("ZZ")

$bookmarkPara = $d.Paragraphs($bookmarkParaIndex)
$anchor = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $anchor)

$placeholder = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start + 2)
$placeholder.Text = ""
